{"js": "// Prepare public release and harden processing reliability:\n// - Retitle the document and simplify the intro paragraph.\n// - Drop the first (Key/Value) table and the \"Here is another table:\" lead-in.\n// - Turn the remaining table into a Name/Age/City roster styled with\n//   \"Light Grid - Accent 1\", trimmed to two sample rows (Alice, Bob).\n\nconst body = context.document.body;\n\n// 1) Update the title and intro paragraph text in place (keeps styles).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"Test Document with Table\", Word.InsertLocation.replace);\nparagraphs.items[1].insertText(\"This is a test document.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Remove the first table (Key/Value, Status/Active) entirely.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\ntables.items[0].delete();\nawait context.sync();\n\n// 3) Remove the now-orphaned \"Here is another table:\" paragraph.\nparagraphs.load(\"items\");\nawait context.sync();\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.text.trim() === \"Here is another table:\") {\n    p.delete();\n  }\n}\nawait context.sync();\n\n// 4) Re-fetch the remaining table (now the only one) and style it.\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.style = \"LightGrid-Accent1\";\nawait context.sync();\n\n// 5) Rewrite header + data cells: Product/Price/Stock -> Name/Age/City,\n//    Apple/1.99/100 -> Alice/30/NYC, Banana/0.99/50 -> Bob/25/LA.\ntable.getCell(0, 0).value = \"Name\";\ntable.getCell(0, 1).value = \"Age\";\ntable.getCell(0, 2).value = \"City\";\n\ntable.getCell(1, 0).value = \"Alice\";\ntable.getCell(1, 1).value = \"30\";\ntable.getCell(1, 2).value = \"NYC\";\n\ntable.getCell(2, 0).value = \"Bob\";\ntable.getCell(2, 1).value = \"25\";\ntable.getCell(2, 2).value = \"LA\";\nawait context.sync();\n\n// 6) Drop the trailing Orange/2.49/75 row.\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nrows.items[rows.items.length - 1].delete();\nawait context.sync();\n", "ps1": "# Prepare public release and harden processing reliability:\n# - Retitle the document and simplify the intro paragraph.\n# - Drop the first (Key/Value) table and the \"Here is another table:\" lead-in.\n# - Turn the remaining table into a Name/Age/City roster styled with\n#   \"Light Grid - Accent 1\", trimmed to two sample rows (Alice, Bob).\n\n$d = $word.ActiveDocument\n\n# 1) Update the title and intro paragraph text in place (keeps styles).\n$d.Paragraphs.Item(1).Range.Text = \"Test Document with Table\"\n$d.Paragraphs.Item(2).Range.Text = \"This is a test document.\"\n\n# 2) Remove the first table (Key/Value, Status/Active) entirely.\n$d.Tables.Item(1).Delete()\n\n# 3) Remove the now-orphaned \"Here is another table:\" paragraph. Re-derive\n#    the paragraph collection from a fresh Range so indices/text reflect\n#    the table deletion above rather than a stale pre-delete snapshot.\n$liveRange = $d.Range()\nforeach ($p in $liveRange.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq \"Here is another table:\") {\n        $p.Range.Delete()\n    }\n}\n\n# 4) Re-fetch the remaining table (now the only one) and style it.\n$table = $d.Tables.Item(1)\n$table.Style = \"LightGrid-Accent1\"\n\n# 5) Rewrite header + data cells: Product/Price/Stock -> Name/Age/City,\n#    Apple/1.99/100 -> Alice/30/NYC, Banana/0.99/50 -> Bob/25/LA.\n$table.Cell(1,1).Range.Text = \"Name\"\n$table.Cell(1,2).Range.Text = \"Age\"\n$table.Cell(1,3).Range.Text = \"City\"\n\n$table.Cell(2,1).Range.Text = \"Alice\"\n$table.Cell(2,2).Range.Text = \"30\"\n$table.Cell(2,3).Range.Text = \"NYC\"\n\n$table.Cell(3,1).Range.Text = \"Bob\"\n$table.Cell(3,2).Range.Text = \"25\"\n$table.Cell(3,3).Range.Text = \"LA\"\n\n# 6) Drop the trailing Orange/2.49/75 row.\n$table = $d.Tables.Item(1)\n$table.Rows.Item($table.Rows.Count).Delete()\n"}
